$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (2..18) for columns:
#   D (Fecha, date serial), I (Calidad), J (Volumen),
#   K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
#   P (Precio $/Kg)
$data = @{}
$data[2]  = @(44253, "Segunda", 1000, 800, 900, 850, 850)
$data[3]  = @(44253, "Tercera", 800, 600, 700, 650, 650)
$data[4]  = @(44229, "Segunda", 760, 550, 600, 575, 575)
$data[5]  = @(44210, "Segunda", 900, 600, 700, 650, 650)
$data[6]  = @(44474, "Segunda", 200, 600, 700, 650, 650)
$data[7]  = @(44174, "Segunda", 800, 450, 500, 475, 475)
$data[8]  = @(44174, "Tercera", 1200, 250, 350, 300, 300)
$data[9]  = @(44573, "Tercera", 800, 600, 650, 625, 625)
$data[10] = @(44245, "Primera", 800, 850, 900, 875, 875)
$data[11] = @(44245, "Segunda", 1000, 750, 800, 775, 775)
$data[12] = @(44544, "Primera", 1000, 600, 650, 625, 625)
$data[13] = @(44658, "Segunda", 1000, 600, 650, 625, 625)
$data[14] = @(44224, "Segunda", 800, 850, 900, 875, 875)
$data[15] = @(44201, "Segunda", 500, 800, 900, 850, 850)
$data[16] = @(44278, "Segunda", 700, 600, 700, 650, 650)
$data[17] = @(44278, "Tercera", 400, 500, 600, 550, 550)
$data[18] = @(44267, "Tercera", 400, 500, 600, 550, 550)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 9).Value  = $vals[1]   # I - Calidad
    $ws.Cells.Item($r, 10).Value = $vals[2]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals[3]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals[4]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals[5]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals[6]   # P - Precio $/Kg
}
